# Weekly update: insert a new price observation as row 172, pushing the
# existing rows 172-255 down to 173-256 (net range grows to A1:R256).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 172..255 down by one, leaving a blank row 172 to populate.
$ws.Rows.Item(172).Insert()

# Populate the newly inserted row 172 with the new weekly observation.
$ws.Range("A172").Value = 4
$ws.Range("B172").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C172").Value = "Los Lagos"
$ws.Range("D172").Value = 44917
$ws.Range("E172").Value = 10
$ws.Range("F172").Value = 100112009
$ws.Range("G172").Value = "Acelga"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 15
$ws.Range("K172").Value = 12000
$ws.Range("L172").Value = 12000
$ws.Range("M172").Value = 12000
$ws.Range("N172").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O172").Value = "Región de La Araucanía"
$ws.Range("P172").Value = 1000
$ws.Range("Q172").Value = 12
$ws.Range("R172").Value = "Hortaliza"
